$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")
$ws.Activate()

# Row 5: Robert's peer assessment (Criterion 1 Online collaboration table)
$ws.Range("B5").Value = "Sufficient"
$ws.Range("C5").Value = "Doing his best at understanding git, participatin in meetings"

# Row 16: Robert's self assessment (Criterion 1 International Collaboration table)
$ws.Range("B16").Value = "Sufficient"
$ws.Range("C16").Value = "Very good comments to discussions some times a bit delayed, " + [char]10
$ws.Range("C16").WrapText = $true

# Leave the view where the author last left it while filling these in
[void]$ws.Range("C18").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
